# Apply the "Task List" sheet updates described in the commit:
#  - MagicIcon / CoinIcon / RangedIcon rows (3-5) move from TODO -> Done,
#    and get actual Time Spent (C) / Over-Under (D) numbers filled in.
#  - IconSpawner / EnemySpawner related rows (6-9) move from TODO -> In Progress.
#  - Selection moves to G13.
#  - The C14 "Time Spent" total recalculates automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tasks 02-04 to 02-11")

# xlPasteFormats
$xlPasteFormats = -4122

# --- Rows 3-5: Create MagicIcon / CoinIcon / RangedIcon : TODO -> Done ---
# Fill in Time Spent (C) and Over/Under (D) values.
$ws.Range("C3").Value2 = 1
$ws.Range("D3").Value2 = 0

$ws.Range("C4").Value2 = 0.5
$ws.Range("D4").Value2 = -0.5

$ws.Range("C5").Value2 = 1.5
$ws.Range("D5").Value2 = 0.5

# Re-use the exact "Done" (green) formatting already present on F2, then set the text.
foreach ($addr in @("F3", "F4", "F5")) {
    $ws.Range("F2").Copy() | Out-Null
    $ws.Range($addr).PasteSpecial($xlPasteFormats) | Out-Null
    $ws.Range($addr).Value2 = "Done"
}

# --- Rows 6-9: IconSpawner / Base Enemy / EnemyPrefab / EnemySpawner : TODO -> In Progress ---
# Re-use the exact "In Progress" (yellow) formatting already present on H9, then set the text.
foreach ($addr in @("F6", "F7", "F8", "F9")) {
    $ws.Range("H9").Copy() | Out-Null
    $ws.Range($addr).PasteSpecial($xlPasteFormats) | Out-Null
    $ws.Range($addr).Value2 = "In Progress"
}

# Clear clipboard marching ants state
$excel.CutCopyMode = 0

# Make sure dependent formulas (C14 = SUM(C3:C13)) are up to date.
$excel.Calculate()

# --- Update the selected cell shown when the workbook is opened ---
$ws.Activate()
$ws.Range("G13").Select() | Out-Null
